# "Delete clase Espia y Add tests comunicadorEspias y fix comEspias y excel"
#
# The spreadsheet used to model a "Espia" (spy) with a name, a paired
# "Compañero" (partner) and a "Probalidad" (probability) in columns A, B
# and C. That data structure is removed: columns B and C are deleted so
# only the plain list of names remains in column A, and a new name
# ("Pepe") is appended to the list. A stray formatted (underlined) cell
# is also left at G8, and the selection / page setup are refreshed to
# match the updated sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "Compañero" / "Probalidad" columns (the old Espia pairing
# data) entirely, collapsing the sheet back down to a single column of
# names.
$ws.Range("B1:C6").ClearContents()
$ws.Columns("B:C").Delete()

# Add the new spy to the list.
$ws.Range("A7").Value = "Pepe"

# Leftover formatting artifact: an empty, underlined cell at G8.
$ws.Range("G8").Font.Underline = $true

# Update the current selection shown when the sheet is opened.
$ws.Range("A10").Select()

# Refresh the page setup for printing.
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
